# "update data with resort sheetname"
# Re-sort the sheet tabs so that "总计" (the summary sheet) comes first,
# ahead of "2022-Q2" (the per-quarter holdings sheet).

$wb = $excel.ActiveWorkbook

# Move "总计" in front of the first sheet so the tab order becomes
# 总计, 2022-Q2.
$summarySheet = $wb.Worksheets.Item("总计")
$summarySheet.Move($wb.Worksheets.Item(1))

# Keep "2022-Q2" as the active/selected sheet, same as before the reorder.
# (Re-fetch the reference after the move so it reflects the new tab order.)
$quarterSheet = $wb.Worksheets.Item("2022-Q2")
$quarterSheet.Activate()
